$wb = $excel.ActiveWorkbook

# --- 1. Reorder tabs first: move "level" to sit right before "type4" ---
$levelTmp = $wb.Worksheets.Item("level")
$type4Tmp = $wb.Worksheets.Item("type4")
$levelTmp.Move($type4Tmp)

# Re-fetch fresh references by name now that positions have changed
# (sheet variables captured before a Move can resolve to a different
# sheet afterwards, since they track position rather than identity).
$enemy = $wb.Worksheets.Item("enemy")
$level = $wb.Worksheets.Item("level")
$type4 = $wb.Worksheets.Item("type4")

# --- 2. Update "enemy" sheet data: adjust two existing values and append a new row ---
$enemy.Range("D5").Value = 5
$enemy.Range("D6").Value = 5
$enemy.Range("A8").Value = 4
$enemy.Range("B8").Value = 4
$enemy.Range("C8").Value = 15
$enemy.Range("D8").Value = 5

# --- 3. Update "level" sheet data (column B index values) ---
$level.Range("B6").Value = 2
$level.Range("B8").Value = 4
$level.Range("B9").Value = 4
$level.Range("B10").Value = 4
$level.Range("B11").Value = 4
$level.Range("B12").Value = 4
$level.Range("B13").Value = 4
$level.Range("B14").Value = 4

# --- 4. Fix up selections on each affected sheet ---
$type4.Activate()
$type4.Range("C20").Select()

$level.Activate()
$level.Range("E16").Select()

# --- 5. Activate "enemy" last so it becomes the active tab/sheet ---
$enemy.Activate()
$enemy.Range("F7").Select()
